$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-04-01 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-02 Tuesday", 2) | Out-Null

# Update the 20x5 practice-problem table
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "7+16="
$t.Cell(1,2).Range.Text = "62-18="
$t.Cell(1,3).Range.Text = "63-34="
$t.Cell(1,4).Range.Text = "40-13="
$t.Cell(1,5).Range.Text = "82-26="

$t.Cell(2,1).Range.Text = "59+19="
$t.Cell(2,2).Range.Text = "96-29="
$t.Cell(2,3).Range.Text = "27+45="
$t.Cell(2,4).Range.Text = "70-28="
$t.Cell(2,5).Range.Text = "8+8="

$t.Cell(3,1).Range.Text = "94-45="
$t.Cell(3,2).Range.Text = "28+15="
$t.Cell(3,3).Range.Text = "22-18="
$t.Cell(3,4).Range.Text = "80-62="
$t.Cell(3,5).Range.Text = "73+8="

$t.Cell(4,1).Range.Text = "49+42="
$t.Cell(4,2).Range.Text = "19+72="
$t.Cell(4,3).Range.Text = "84+9="
$t.Cell(4,4).Range.Text = "51-48="
$t.Cell(4,5).Range.Text = "45+17="

$t.Cell(5,1).Range.Text = "28+68="
$t.Cell(5,2).Range.Text = "29+42="
$t.Cell(5,3).Range.Text = "15-9="
$t.Cell(5,4).Range.Text = "36+58="
$t.Cell(5,5).Range.Text = "25+29="

$t.Cell(6,1).Range.Text = "69+2="
$t.Cell(6,2).Range.Text = "20-2="
$t.Cell(6,3).Range.Text = "54-35="
$t.Cell(6,4).Range.Text = "43-37="
$t.Cell(6,5).Range.Text = "30-24="

$t.Cell(7,1).Range.Text = "46-39="
$t.Cell(7,2).Range.Text = "5+18="
$t.Cell(7,3).Range.Text = "82-17="
$t.Cell(7,4).Range.Text = "80-29="
$t.Cell(7,5).Range.Text = "23-5="

$t.Cell(8,1).Range.Text = "5+39="
$t.Cell(8,2).Range.Text = "86-8="
$t.Cell(8,3).Range.Text = "82-43="
$t.Cell(8,4).Range.Text = "47+18="
$t.Cell(8,5).Range.Text = "52+39="

$t.Cell(9,1).Range.Text = "70-54="
$t.Cell(9,2).Range.Text = "27+15="
$t.Cell(9,3).Range.Text = "82-57="
$t.Cell(9,4).Range.Text = "32-24="
$t.Cell(9,5).Range.Text = "50-29="

$t.Cell(10,1).Range.Text = "48+14="
$t.Cell(10,2).Range.Text = "34+49="
$t.Cell(10,3).Range.Text = "17+6="
$t.Cell(10,4).Range.Text = "94-56="
$t.Cell(10,5).Range.Text = "17+78="

$t.Cell(11,1).Range.Text = "69+2="
$t.Cell(11,2).Range.Text = "43-17="
$t.Cell(11,3).Range.Text = "8+49="
$t.Cell(11,4).Range.Text = "63-8="
$t.Cell(11,5).Range.Text = "27+6="

$t.Cell(12,1).Range.Text = "55-8="
$t.Cell(12,2).Range.Text = "46+35="
$t.Cell(12,3).Range.Text = "19+18="
$t.Cell(12,4).Range.Text = "49+35="
$t.Cell(12,5).Range.Text = "81-17="

$t.Cell(13,1).Range.Text = "15+9="
$t.Cell(13,2).Range.Text = "14-8="
$t.Cell(13,3).Range.Text = "78+19="
$t.Cell(13,4).Range.Text = "36+16="
$t.Cell(13,5).Range.Text = "80-31="

$t.Cell(14,1).Range.Text = "41-29="
$t.Cell(14,2).Range.Text = "81-58="
$t.Cell(14,3).Range.Text = "62-26="
$t.Cell(14,4).Range.Text = "23+69="
$t.Cell(14,5).Range.Text = "79+5="

$t.Cell(15,1).Range.Text = "8+87="
$t.Cell(15,2).Range.Text = "47+7="
$t.Cell(15,3).Range.Text = "29+29="
$t.Cell(15,4).Range.Text = "86-59="
$t.Cell(15,5).Range.Text = "7+88="

$t.Cell(16,1).Range.Text = "19+73="
$t.Cell(16,2).Range.Text = "44+49="
$t.Cell(16,3).Range.Text = "94-49="
$t.Cell(16,4).Range.Text = "24+47="
$t.Cell(16,5).Range.Text = "47+16="

$t.Cell(17,1).Range.Text = "55-36="
$t.Cell(17,2).Range.Text = "52-16="
$t.Cell(17,3).Range.Text = "5+37="
$t.Cell(17,4).Range.Text = "21-12="
$t.Cell(17,5).Range.Text = "18+33="

$t.Cell(18,1).Range.Text = "28+64="
$t.Cell(18,2).Range.Text = "40-6="
$t.Cell(18,3).Range.Text = "59+32="
$t.Cell(18,4).Range.Text = "75-27="
$t.Cell(18,5).Range.Text = "37+6="

$t.Cell(19,1).Range.Text = "97-19="
$t.Cell(19,2).Range.Text = "31-13="
$t.Cell(19,3).Range.Text = "81-75="
$t.Cell(19,4).Range.Text = "29+38="
$t.Cell(19,5).Range.Text = "37+37="

$t.Cell(20,1).Range.Text = "18-9="
$t.Cell(20,2).Range.Text = "38+19="
$t.Cell(20,3).Range.Text = "91-46="
$t.Cell(20,4).Range.Text = "79+7="
$t.Cell(20,5).Range.Text = "18+5="
